# Updates cryptos list price/volume columns (Sat Jun 15 20:28:10 UTC 2024 run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# Row 2
Set-TextCell "D2" "65.980.71"
$ws.Range("E2").Value = "  +0.35%  "

# Row 3
Set-TextCell "D3" "3.545.91"
$ws.Range("E3").Value = "  +3.79%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
Set-TextCell "D5" "605.62"
$ws.Range("E5").Value = "  +1.47%  "

# Row 6
Set-TextCell "D6" "144.48"
$ws.Range("E6").Value = "  +1.38%  "

# Row 7
Set-TextCell "D7" "3.543.17"
$ws.Range("E7").Value = "  +3.76%  "

# Row 8
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
Set-TextCell "D9" "0.492"
$ws.Range("E9").Value = "  +4.63%  "

# Row 10
Set-TextCell "D10" "0.136"
$ws.Range("E10").Value = "  +1.56%  "

# Row 11
Set-TextCell "D11" "7.93"
$ws.Range("E11").Value = "  -0.76%  "

# Row 12
Set-TextCell "D12" "0.412"
$ws.Range("E12").Value = "  +1.75%  "

# Row 13
Set-TextCell "D13" "4.146.54"
$ws.Range("E13").Value = "  +3.90%  "

# Row 14
Set-TextCell "D14" "0.0000206"
$ws.Range("E14").Value = "  +2.89%  "

# Row 15
Set-TextCell "D15" "30.00"
$ws.Range("E15").Value = "  +1.30%  "

# Row 16
Set-TextCell "D16" "3.545.85"
$ws.Range("E16").Value = "  +3.92%  "

# Row 17
Set-TextCell "D17" "66.189.33"
$ws.Range("E17").Value = "  +0.78%  "

# Row 18
$ws.Range("E18").Value = "  -0.52%  "

# Row 19
Set-TextCell "D19" "11.29"
$ws.Range("E19").Value = "  +8.69%  "

# Row 20
Set-TextCell "D20" "6.17"
$ws.Range("E20").Value = "  +0.84%  "

# Row 21
Set-TextCell "D21" "14.84"
$ws.Range("E21").Value = "  +1.70%  "

# Row 22
Set-TextCell "D22" "429.80"
$ws.Range("E22").Value = "  +3.49%  "

# Row 23
Set-TextCell "D23" "0.609"
$ws.Range("E23").Value = "  +5.44%  "

# Row 24
Set-TextCell "D24" "79.12"
$ws.Range("E24").Value = "  +2.43%  "

# Row 25
Set-TextCell "D25" "3.688.26"
$ws.Range("E25").Value = "  +3.88%  "

# Row 26
$ws.Range("E26").Value = "  +0.04%  "

# Row 27
Set-TextCell "D27" "0.0000117"
$ws.Range("E27").Value = "  +7.73%  "

# Row 28
Set-TextCell "D28" "2.51"
$ws.Range("E28").Value = "  +3.87%  "

# Row 29
Set-TextCell "D29" "7.94"
$ws.Range("E29").Value = "  +1.34%  "

# Row 30
Set-TextCell "D30" "9.05"
$ws.Range("E30").Value = "  -2.16%  "

# Row 31
Set-TextCell "D31" "1.00"
$ws.Range("E31").Value = "  +0.11%  "

# Row 32
Set-TextCell "D32" "1.46"
$ws.Range("E32").Value = "  +0.50%  "

# Row 33
Set-TextCell "D33" "25.47"
$ws.Range("E33").Value = "  +3.76%  "

# Row 34
Set-TextCell "D34" "3.542.17"
$ws.Range("E34").Value = "  +3.84%  "

# Row 35
Set-TextCell "D35" "0.153"
$ws.Range("E35").Value = "  -4.06%  "

# Row 36
$ws.Range("E36").Value = "  +0.08%  "

# Row 37
Set-TextCell "D37" "1.74"
$ws.Range("E37").Value = "  +3.09%  "

# Row 38
Set-TextCell "D38" "7.87"
$ws.Range("E38").Value = "  +4.59%  "

# Row 39
Set-TextCell "D39" "5.59"
$ws.Range("E39").Value = "  +1.73%  "

# Row 40
Set-TextCell "D40" "0.999"
$ws.Range("E40").Value = "  +0.07%  "

# Row 41
Set-TextCell "D41" "174.78"
$ws.Range("E41").Value = "  +3.70%  "

# Row 42
Set-TextCell "D42" "0.0849"
$ws.Range("E42").Value = "  -0.64%  "

# Row 43
Set-TextCell "D43" "5.19"
$ws.Range("E43").Value = "  +3.15%  "

# Row 44
Set-TextCell "D44" "0.891"
$ws.Range("E44").Value = "  +2.07%  "

# Row 45
$ws.Range("E45").Value = "  +0.99%  "

# Row 46
Set-TextCell "D46" "45.99"
$ws.Range("E46").Value = "  +1.32%  "

# Row 47
Set-TextCell "D47" "1.20"
$ws.Range("E47").Value = "  +2.12%  "

# Row 48
Set-TextCell "D48" "25.39"
$ws.Range("E48").Value = "  -3.08%  "

# Rows 49-51: reorder coins (dogwifhat/Cosmos/EnergySwap) with updated price/volume
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D49" "7.11"
$ws.Range("E49").Value = "  +0.89%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D50" "23.34"
$ws.Range("E50").Value = "  +13.58%  "

$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell "D51" "2.33"
$ws.Range("E51").Value = "  +2.79%  "
